$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -7
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -7
$ws.Range("F10").Value = 8
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 8
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = 6
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = 3
$ws.Range("F20").Value = -3
$ws.Range("F21").Value = -3
